$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 8
$ws1.Range("F6").Value = 524
$ws1.Range("F7").Value = 1574
$ws1.Range("F9").Value = 10
$ws1.Range("F10").Value = 1276
$ws1.Range("F13").Value = 169
$ws1.Range("F15").Value = 6
$ws1.Range("F16").Value = 6
$ws1.Range("F19").Value = 194
$ws1.Range("F20").Value = 186

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 8
$ws4.Range("F6").Value = 524
$ws4.Range("F7").Value = 1574
$ws4.Range("F9").Value = 6
$ws4.Range("F10").Value = 10
$ws4.Range("F11").Value = 1276
$ws4.Range("F14").Value = 169
$ws4.Range("F16").Value = 6
$ws4.Range("F17").Value = 6
$ws4.Range("F20").Value = 194
$ws4.Range("F21").Value = 186

$wb.Save()
